$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 321 (shifts existing rows 321:407 down to 322:408),
# matching the new weekly price record that was added to the dataset.
$ws.Rows.Item(321).Insert()

# Populate the newly inserted row 321 with the new weekly record.
$ws.Range("A321").Value = 8
$ws.Range("B321").Value = "Terminal La Palmera de La Serena"
$ws.Range("C321").Value = "Coquimbo"
$ws.Range("D321").Value = 44943
$ws.Range("E321").Value = 4
$ws.Range("F321").Value = 100112003
$ws.Range("G321").Value = "Ajo"
$ws.Range("H321").Value = "Chino"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 400
$ws.Range("K321").Value = 16000
$ws.Range("L321").Value = 17000
$ws.Range("M321").Value = 16500
$ws.Range("N321").Value = '$/caja 10 kilos'
$ws.Range("O321").Value = "China"
$ws.Range("P321").Value = 1650
$ws.Range("Q321").Value = 10
$ws.Range("R321").Value = "Hortaliza"

# Make sure the date column keeps the same date/time number format used by the rest of column D.
$ws.Range("D321").NumberFormat = $ws.Range("D322").NumberFormat
